$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data,
# keeping values stored as plain text (matching the original inline-string cells).
$cellUpdates = @{
    'D2' = '303.67'
    'E2' = '5.67%'
    'D3' = '32.16'
    'E3' = '10.39%'
    'E4' = '1.43%'
    'D5' = '0.07480'
    'E5' = '7.49%'
    'D6' = '7.852'
    'D7' = '3.811'
    'D8' = '1.525'
    'E8' = '9.11%'
    'D9' = '0.9191'
    'E9' = '2.19%'
    'D10' = '0.1685'
    'E10' = '5.11%'
    'D11' = '0.08026'
    'E11' = '4.10%'
    'D12' = '0.08076'
    'E12' = '4.41%'
    'D13' = '0.03008'
    'E13' = '2.44%'
    'E14' = '9.84%'
    'D15' = '0.001497'
    'E15' = '-5.92%'
    'D16' = '0.04599'
    'E16' = '1.57%'
    'D17' = '0.006451'
    'E17' = '-1.13%'
    'D18' = '3.467'
    'E18' = '-0.51%'
    'D19' = '2.232'
    'E19' = '0.08%'
    'E20' = '2.11%'
    'D21' = '0.1342'
    'E21' = '0.36%'
    'D22' = '4.475'
    'E22' = '10.72%'
    'D24' = '0.001216'
    'E24' = '0.42%'
    'D25' = '0.004450'
    'E25' = '7.49%'
    'D26' = '0.0001399'
    'E26' = '19.50%'
    'E27' = '6.33%'
    'D39' = '0.01722'
    'E39' = '2,548.11%'
    'D40' = '0.04505'
    'E40' = '3.43%'
    'D41' = '0.007162'
    'E41' = '3.23%'
    'D42' = '0.1348'
    'E42' = '8.21%'
    'D43' = '0.002248'
    'E43' = '8.55%'
    'D44' = '0.01280'
    'E44' = '9.14%'
    'D45' = '0.00006167'
    'E45' = '5.68%'
    'D46' = '1.868'
    'E46' = '-3.16%'
    'D47' = '0.01299'
    'E47' = '-0.10%'
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
